# "added forced ownership over send requests"
# The NORM ("G") column was sharing the same divisor (1) for every benchmark,
# which made the normalized ratios meaningless. Force the NORM column to take
# ownership of (i.e. use) the same per-benchmark "Shared" request counts as
# column F, rather than the placeholder value of 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the "Shared" request counts (column F) into the "NORM" column (G)
# for each benchmark row so G actually reflects real request ownership.
$ws.Range("G3").Value = $ws.Range("F3").Value()
$ws.Range("G4").Value = $ws.Range("F4").Value()
$ws.Range("G5").Value = $ws.Range("F5").Value()
$ws.Range("G6").Value = $ws.Range("F6").Value()

# Reflect the author's last cursor/selection position before saving.
$ws.Range("H19").Select()

$wb.Save()
